# stationData.xlsx update:
#  - Add CONNECTIONS entries for Aberdeen, Braid, Brentwood Town Center
#    and Bridgeport (rows 4-7, column G)
#  - Fix Sapperton's station CODE from "SPT" to "SAP"
#  - Leave the active selection on G8, matching the author's last edit

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G4").Value = "BRP-CAN-WFRONT:LDN-CAN-RICHBR"
$ws.Range("B43").Value = "SAP"
$ws.Range("G5").Value = "SAP-EXPO-WFRONT:LTC-EXPO-PWAYU"
$ws.Range("G6").Value = "HLD-MILL-LLDOUG:GLM-MILL-VCCCL"
$ws.Range("G7").Value = "MRD-CAN-WFRONT:ABD-CAN-RICHBR:TPL-CAN-YVRA"

$ws.Range("G8").Select()
